$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.816965666666666
$ws.Range("H2").Value = 17.450897
$ws.Range("I2").Value = 0.03136298918947451
$ws.Range("J2").Value = 0.03338421840610672
$ws.Range("M2").Value = 16.92841533333333
$ws.Range("N2").Value = 50.785246
$ws.Range("O2").Value = 0.0242040298661412
$ws.Range("P2").Value = 0.0261568218313686
$ws.Range("Q2").Value = 98.47201078507354
$ws.Range("R2").Value = 886.2480970656619
$ws.Range("S2").Value = 0.0007591107270335047
$ws.Range("T2").Value = 0.0008732250528280298
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.816965666666666
$ws.Range("H3").Value = 17.450897
$ws.Range("I3").Value = 0.03136298918947451
$ws.Range("J3").Value = 0.03338421840610672
$ws.Range("O3").Value = 0.2723327394629209
$ws.Range("P3").Value = 0.2943046668003394
$ws.Range("Q3").Value = 1107.962294123412
$ws.Range("R3").Value = 9971.660647110704
$ws.Range("S3").Value = 0.008541168763715568
$ws.Range("T3").Value = 0.009825131274398996
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.816965666666666
$ws.Range("H4").Value = 17.450897
$ws.Range("I4").Value = 0.03136298918947451
$ws.Range("J4").Value = 0.03338421840610672
$ws.Range("M4").Value = 110.8005546666667
$ws.Range("N4").Value = 332.401664
$ws.Range("O4").Value = 0.1584212037293475
$ws.Range("P4").Value = 0.17120269736802
$ws.Range("Q4").Value = 644.523022343623
$ws.Range("R4").Value = 5800.707201092607
$ws.Range("S4").Value = 0.004968562499947066
$ws.Range("T4").Value = 0.005715468240648573
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.816965666666666
$ws.Range("H5").Value = 17.450897
$ws.Range("I5").Value = 0.03136298918947451
$ws.Range("J5").Value = 0.03338421840610672
$ws.Range("M5").Value = 156.6465685
$ws.Range("N5").Value = 313.293137
$ws.Range("O5").Value = 0.223971242892229
$ws.Range("P5").Value = 0.1613608953572767
$ws.Range("Q5").Value = 911.207710765648
$ws.Range("R5").Value = 5467.246264593889
$ws.Range("S5").Value = 0.007024407669582147
$ws.Range("T5").Value = 0.005386907372812257
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.816965666666666
$ws.Range("H6").Value = 17.450897
$ws.Range("I6").Value = 0.03136298918947451
$ws.Range("J6").Value = 0.03338421840610672
$ws.Range("M6").Value = 224.5584563333333
$ws.Range("N6").Value = 673.675369
$ws.Range("O6").Value = 0.3210707840493613
$ws.Range("P6").Value = 0.3469749186429952
$ws.Range("Q6").Value = 1306.248830650666
$ws.Range("R6").Value = 11756.23947585599
$ws.Range("S6").Value = 0.01006973952919623
$ws.Range("T6").Value = 0.01158348646541886
$ws.Range("I7").Value = 0.7673644067452855
$ws.Range("J7").Value = 0.8168182183493699
$ws.Range("M7").Value = 16.92841533333333
$ws.Range("N7").Value = 50.785246
$ws.Range("O7").Value = 0.0242040298661412
$ws.Range("P7").Value = 0.0261568218313686
$ws.Range("Q7").Value = 2409.333998127408
$ws.Range("R7").Value = 21684.00598314668
$ws.Range("S7").Value = 0.01857331101907661
$ws.Range("T7").Value = 0.0213653686059804
$ws.Range("I8").Value = 0.7673644067452855
$ws.Range("J8").Value = 0.8168182183493699
$ws.Range("O8").Value = 0.2723327394629209
$ws.Range("P8").Value = 0.2943046668003394
$ws.Range("S8").Value = 0.2089784510552827
$ws.Range("T8").Value = 0.2403934135877582
$ws.Range("I9").Value = 0.7673644067452855
$ws.Range("J9").Value = 0.8168182183493699
$ws.Range("M9").Value = 110.8005546666667
$ws.Range("N9").Value = 332.401664
$ws.Range("O9").Value = 0.1584212037293475
$ws.Range("P9").Value = 0.17120269736802
$ws.Range("Q9").Value = 15769.67117791107
$ws.Range("R9").Value = 141927.0406011996
$ws.Range("S9").Value = 0.1215667930156448
$ws.Range("T9").Value = 0.1398414822407525
$ws.Range("I10").Value = 0.7673644067452855
$ws.Range("J10").Value = 0.8168182183493699
$ws.Range("M10").Value = 156.6465685
$ws.Range("N10").Value = 313.293137
$ws.Range("O10").Value = 0.223971242892229
$ws.Range("P10").Value = 0.1613608953572767
$ws.Range("Q10").Value = 22294.69774609602
$ws.Range("R10").Value = 133768.1864765761
$ws.Range("S10").Value = 0.1718675599299995
$ws.Range("T10").Value = 0.1318025190569898
$ws.Range("I11").Value = 0.7673644067452855
$ws.Range("J11").Value = 0.8168182183493699
$ws.Range("M11").Value = 224.5584563333333
$ws.Range("N11").Value = 673.675369
$ws.Range("O11").Value = 0.3210707840493613
$ws.Range("P11").Value = 0.3469749186429952
$ws.Range("Q11").Value = 31960.24629343583
$ws.Range("R11").Value = 287642.2166409225
$ws.Range("S11").Value = 0.2463782917252818
$ws.Range("T11").Value = 0.2834154348578889
$ws.Range("G12").Value = 1.295824333333333
$ws.Range("H12").Value = 3.887473
$ws.Range("I12").Value = 0.006986619293745993
$ws.Range("J12").Value = 0.007436881191828873
$ws.Range("M12").Value = 16.92841533333333
$ws.Range("N12").Value = 50.785246
$ws.Range("O12").Value = 0.0242040298661412
$ws.Range("P12").Value = 0.0261568218313686
$ws.Range("Q12").Value = 21.93625251370644
$ws.Range("R12").Value = 197.426272623358
$ws.Range("S12").Value = 0.0001691043420491863
$ws.Range("T12").Value = 0.000194525176315724
$ws.Range("G13").Value = 1.295824333333333
$ws.Range("H13").Value = 3.887473
$ws.Range("I13").Value = 0.006986619293745993
$ws.Range("J13").Value = 0.007436881191828873
$ws.Range("O13").Value = 0.2723327394629209
$ws.Range("P13").Value = 0.2943046668003394
$ws.Range("Q13").Value = 246.8167397597282
$ws.Range("R13").Value = 2221.350657837554
$ws.Range("S13").Value = 0.001902685171850344
$ws.Range("T13").Value = 0.002188708841194908
$ws.Range("G14").Value = 1.295824333333333
$ws.Range("H14").Value = 3.887473
$ws.Range("I14").Value = 0.006986619293745993
$ws.Range("J14").Value = 0.007436881191828873
$ws.Range("M14").Value = 110.8005546666667
$ws.Range("N14").Value = 332.401664
$ws.Range("O14").Value = 0.1584212037293475
$ws.Range("P14").Value = 0.17120269736802
$ws.Range("Q14").Value = 143.5780548838969
$ws.Range("R14").Value = 1292.202493955072
$ws.Range("S14").Value = 0.001106828638513924
$ws.Range("T14").Value = 0.001273214120046599
$ws.Range("G15").Value = 1.295824333333333
$ws.Range("H15").Value = 3.887473
$ws.Range("I15").Value = 0.006986619293745993
$ws.Range("J15").Value = 0.007436881191828873
$ws.Range("M15").Value = 156.6465685
$ws.Range("N15").Value = 313.293137
$ws.Range("O15").Value = 0.223971242892229
$ws.Range("P15").Value = 0.1613608953572767
$ws.Range("Q15").Value = 202.9864351954668
$ws.Range("R15").Value = 1217.918611172801
$ws.Range("S15").Value = 0.001564801806835117
$ws.Range("T15").Value = 0.001200021807779198
$ws.Range("G16").Value = 1.295824333333333
$ws.Range("H16").Value = 3.887473
$ws.Range("I16").Value = 0.006986619293745993
$ws.Range("J16").Value = 0.007436881191828873
$ws.Range("M16").Value = 224.5584563333333
$ws.Range("N16").Value = 673.675369
$ws.Range("O16").Value = 0.3210707840493613
$ws.Range("P16").Value = 0.3469749186429952
$ws.Range("Q16").Value = 290.9883119725041
$ws.Range("R16").Value = 2618.894807752537
$ws.Range("S16").Value = 0.002243199334497421
$ws.Range("T16").Value = 0.002580411246492445
$ws.Range("G17").Value = 33.687956
$ws.Range("H17").Value = 67.375912
$ws.Range("I17").Value = 0.181633356699686
$ws.Range("J17").Value = 0.1288926386717328
$ws.Range("M17").Value = 16.92841533333333
$ws.Range("N17").Value = 50.785246
$ws.Range("O17").Value = 0.0242040298661412
$ws.Range("P17").Value = 0.0261568218313686
$ws.Range("Q17").Value = 570.2837108990586
$ws.Range("R17").Value = 3421.702265394352
$ws.Range("S17").Value = 0.004396259190246677
$ws.Range("T17").Value = 0.003371421785111487
$ws.Range("G18").Value = 33.687956
$ws.Range("H18").Value = 67.375912
$ws.Range("I18").Value = 0.181633356699686
$ws.Range("J18").Value = 0.1288926386717328
$ws.Range("O18").Value = 0.2723327394629209
$ws.Range("P18").Value = 0.2943046668003394
$ws.Range("Q18").Value = 6416.57303015803
$ws.Range("R18").Value = 38499.43818094818
$ws.Range("S18").Value = 0.04946470960787136
$ws.Range("T18").Value = 0.03793370507730087
$ws.Range("G19").Value = 33.687956
$ws.Range("H19").Value = 67.375912
$ws.Range("I19").Value = 0.181633356699686
$ws.Range("J19").Value = 0.1288926386717328
$ws.Range("M19").Value = 110.8005546666667
$ws.Range("N19").Value = 332.401664
$ws.Range("O19").Value = 0.1584212037293475
$ws.Range("P19").Value = 0.17120269736802
$ws.Range("Q19").Value = 3732.644210386261
$ws.Range("R19").Value = 22395.86526231757
$ws.Range("S19").Value = 0.0287745750057662
$ws.Range("T19").Value = 0.02206676741148223
$ws.Range("G20").Value = 33.687956
$ws.Range("H20").Value = 67.375912
$ws.Range("I20").Value = 0.181633356699686
$ws.Range("J20").Value = 0.1288926386717328
$ws.Range("M20").Value = 156.6465685
$ws.Range("N20").Value = 313.293137
$ws.Range("O20").Value = 0.223971242892229
$ws.Range("P20").Value = 0.1613608953572767
$ws.Range("Q20").Value = 5277.102707178986
$ws.Range("R20").Value = 21108.41082871594
$ws.Range("S20").Value = 0.04068064865071623
$ws.Range("T20").Value = 0.02079823158103276
$ws.Range("G21").Value = 33.687956
$ws.Range("H21").Value = 67.375912
$ws.Range("I21").Value = 0.181633356699686
$ws.Range("J21").Value = 0.1288926386717328
$ws.Range("M21").Value = 224.5584563333333
$ws.Range("N21").Value = 673.675369
$ws.Range("O21").Value = 0.3210707840493613
$ws.Range("P21").Value = 0.3469749186429952
$ws.Range("Q21").Value = 7564.915396385255
$ws.Range("R21").Value = 45389.49237831153
$ws.Range("S21").Value = 0.05831716424508549
$ws.Range("T21").Value = 0.04472251281680548
$ws.Range("G22").Value = 2.346712
$ws.Range("H22").Value = 7.040136
$ws.Range("I22").Value = 0.01265262807180802
$ws.Range("J22").Value = 0.01346804338096171
$ws.Range("M22").Value = 16.92841533333333
$ws.Range("N22").Value = 50.785246
$ws.Range("O22").Value = 0.0242040298661412
$ws.Range("P22").Value = 0.0261568218313686
$ws.Range("Q22").Value = 39.72611540371734
$ws.Range("R22").Value = 357.535038633456
$ws.Range("S22").Value = 0.0003062445877352179
$ws.Range("T22").Value = 0.0003522812111329587
$ws.Range("G23").Value = 2.346712
$ws.Range("H23").Value = 7.040136
$ws.Range("I23").Value = 0.01265262807180802
$ws.Range("J23").Value = 0.01346804338096171
$ws.Range("O23").Value = 0.2723327394629209
$ws.Range("P23").Value = 0.2943046668003394
$ws.Range("Q23").Value = 446.9801886688587
$ws.Range("R23").Value = 4022.821698019729
$ws.Range("S23").Value = 0.003445724864200933
$ws.Range("T23").Value = 0.003963708019686453
$ws.Range("G24").Value = 2.346712
$ws.Range("H24").Value = 7.040136
$ws.Range("I24").Value = 0.01265262807180802
$ws.Range("J24").Value = 0.01346804338096171
$ws.Range("M24").Value = 110.8005546666667
$ws.Range("N24").Value = 332.401664
$ws.Range("O24").Value = 0.1584212037293475
$ws.Range("P24").Value = 0.17120269736802
$ws.Range("Q24").Value = 260.0169912429226
$ws.Range("R24").Value = 2340.152921186304
$ws.Range("S24").Value = 0.00200444456947556
$ws.Range("T24").Value = 0.002305765355090153
$ws.Range("G25").Value = 2.346712
$ws.Range("H25").Value = 7.040136
$ws.Range("I25").Value = 0.01265262807180802
$ws.Range("J25").Value = 0.01346804338096171
$ws.Range("M25").Value = 156.6465685
$ws.Range("N25").Value = 313.293137
$ws.Range("O25").Value = 0.223971242892229
$ws.Range("P25").Value = 0.1613608953572767
$ws.Range("Q25").Value = 367.604382057772
$ws.Range("R25").Value = 2205.626292346632
$ws.Range("S25").Value = 0.002833824835095949
$ws.Range("T25").Value = 0.002173215538662625
$ws.Range("G26").Value = 2.346712
$ws.Range("H26").Value = 7.040136
$ws.Range("I26").Value = 0.01265262807180802
$ws.Range("J26").Value = 0.01346804338096171
$ws.Range("M26").Value = 224.5584563333333
$ws.Range("N26").Value = 673.675369
$ws.Range("O26").Value = 0.3210707840493613
$ws.Range("P26").Value = 0.3469749186429952
$ws.Range("Q26").Value = 526.9740241789094
$ws.Range("R26").Value = 4742.766217610185
$ws.Range("S26").Value = 0.00406238921530036
$ws.Range("T26").Value = 0.00467307325638952
